$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.919.94'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '1.670.88'
$ws.Range("E3").Value = '  +1.21%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.517'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.95%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0619'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0889'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.45%  '
$ws.Range("D12").Value = '1.906.83'
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("D13").Value = '1.673.94'
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("E14").Value = '  +0.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.525'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("D17").Value = '26.920.87'
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.01'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '232.99'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.21%  '
$ws.Range("D20").Value = '0.0₃0732'
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("E23").Value = '  -1.91%  '
$ws.Range("E24").Value = '  -1.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.94'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.90'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.65%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("E29").Value = '  -1.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0497'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("D33").Value = '1.451.62'
$ws.Range("E33").Value = '  -6.77%  '
$ws.Range("E34").Value = '  +1.33%  '
$ws.Range("E35").Value = '  +1.18%  '
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.578'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.898'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.55%  '
$ws.Range("E39").Value = '  +0.54%  '
$ws.Range("E40").Value = '  +13.22%  '
$ws.Range("E41").Value = '  -4.16%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.23'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("D45").Value = '1.811.44'
$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.780'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.37'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("E48").Value = '  +0.96%  '
$ws.Range("E49").Value = '  +2.24%  '
$ws.Range("E50").Value = '  +0.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.63'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.26%  '
